$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) row above the current row 2 - this shifts the old
# row 2 ("SSO" usage entry) and row 3 ("test2..test5" entry) down to rows
# 3 and 4, and grows the used range to A1:E4.
$ws.Rows.Item(2).Insert()

# Centre-align the (now two-row-tall) header band. Doing the hyperlinked
# A column separately from the rest keeps the hyperlink font on its own
# style record, same as the rest of the plain header cells.
$ws.Range("A1:A2").HorizontalAlignment = -4108
$ws.Range("B1:E2").HorizontalAlignment = -4108

# Merge each header column across the title row + new blank row beneath it.
$ws.Range("A1:A2").Merge()
$ws.Range("B1:B2").Merge()
$ws.Range("C1:C2").Merge()
$ws.Range("D1:D2").Merge()
$ws.Range("E1:E2").Merge()

# The row insert shifted the data rows down but left the hyperlink
# anchors pointing at their old (now wrong) cells. Rebuild the two
# hyperlinks that moved (A2 -> A3, A3 -> A4); the one on A1 never moved.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:supppyy@hotmail.com", "", "", "supppyy@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:supppyy@hotmail.com", "", "", "supppyy@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:supppyy@hotmail.com", "", "", "supppyy@hotmail.com")

# Match the author's final cursor position.
$ws.Range("N27").Select()
